$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D:E").Insert()
$ws.Range("F7:G102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 57100
$ws.Range("E8").Value = 50700
$ws.Range("D9").Value = 20100
$ws.Range("E9").Value = 15600
$ws.Range("D10").Value = 37000
$ws.Range("E10").Value = 35100
$ws.Range("D12").Value = 3500
$ws.Range("E12").Value = 4700
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value = 8700
$ws.Range("E15").Value = 8500
$ws.Range("D17").Value = 45700
$ws.Range("E17").Value = 40600
$ws.Range("D18").Value = 11400
$ws.Range("E18").Value = 10100
$ws.Range("D20").Value = -500
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = 19700
$ws.Range("E21").Value = 18700
$ws.Range("D22").Value = 3600
$ws.Range("E22").Value = 3800
$ws.Range("D23").Value = 7300
$ws.Range("E23").Value = 6400
$ws.Range("D24").Value = 1900
$ws.Range("E24").Value = 1300
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 5400
$ws.Range("E26").Value = 5000
$ws.Range("D27").Value = 5400
$ws.Range("E27").Value = 5000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 500
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 5400
$ws.Range("E33").Value = 5000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 5400
$ws.Range("E35").Value = 5000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 43000
$ws.Range("E41").Value = 44100
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 64800
$ws.Range("E43").Value = 67600
$ws.Range("D44").Value = 40500
$ws.Range("E44").Value = 40000
$ws.Range("D45").Value = 4500
$ws.Range("E45").Value = 5000
$ws.Range("D46").Value = 152900
$ws.Range("E46").Value = 156800
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 38100
$ws.Range("E48").Value = 37400
$ws.Range("D49").Value = 205200
$ws.Range("E49").Value = 213700
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 34500
$ws.Range("E52").Value = 31500
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 430600
$ws.Range("E54").Value = 439400
$ws.Range("D57").Value = 8900
$ws.Range("E57").Value = 7300
$ws.Range("D58").Value = 115700
$ws.Range("E58").Value = 5700
$ws.Range("D59").Value = 40900
$ws.Range("E59").Value = 34000
$ws.Range("D60").Value = 165500
$ws.Range("E60").Value = 46900
$ws.Range("D61").Value = 67300
$ws.Range("E61").Value = 200100
$ws.Range("D62").Value = 500
$ws.Range("E62").Value = "NA"
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 233300
$ws.Range("E66").Value = 247000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 11500
$ws.Range("E72").Value = 6100
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 197300
$ws.Range("E76").Value = 192500
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 5400
$ws.Range("E81").Value = 5000
$ws.Range("D83").Value = 8700
$ws.Range("E83").Value = 8500
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 27200
$ws.Range("E89").Value = 8300
$ws.Range("D91").Value = -1000
$ws.Range("E91").Value = -1400
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -400
$ws.Range("E94").Value = -18400
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -27900
$ws.Range("E100").Value = -700
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -1100
$ws.Range("E102").Value = -10900
